$d = $word.ActiveDocument

# --- Edit 1: replace the two placeholder "Lp <3" bullet paragraphs (items 1 and 2)
# with the real answers to questions a) and b). This also removes the spell-check
# proofErr wrappers that surrounded "Lp" and splits the answer to b) into two
# paragraphs (second one un-bulleted, continuing the same list item visually).
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$rng1 = $d.Range($p1.Range.Start, $p2.Range.End)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00001787" w:rsidRDefault="004926C7" w:rsidP="004926C7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>L</w:t></w:r><w:r><w:t>e temps de préemption représente le temps durant lequel une tâche est prête mais n’est pas exécutée au profit d’une autre tâche dont la priorité est plus élevée. Le temps de blocage est le temps durant lequel une tâche n’est pas exécutée au profit d’une autre tâche dont la priorité est moins élevée, parce que la tâche prioritaire attend après l’autre, généralement parce qu’elle a besoin d’une ressource qui est déjà utilisée par la tâche moins prioritaire.</w:t></w:r></w:p>
<w:p w:rsidR="004926C7" w:rsidRDefault="004926C7" w:rsidP="004926C7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Le protocole ICCP demande peu de changement de contexte, parce que lorsqu’une tâche nécessite une ressource prioritaire, elle acquiert directement une plus grande priorité. Alors que dans le cas de l’héritage de priorité, si une tâche plus prioritaire nécessitant la ressource démarre, il y aura changement de contexte au démarrage de la tâche prioritaire, puis encore une fois lorsque la tâche prioritaire devient bloquée par la ressource et une troisième fois lorsque la tâche moins prioritaire termine d’utiliser la ressource. </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Aussi, l’ICCP permet d’éviter des problèmes d’inter-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bloquage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> qui pourraient survenir dans le cas de l’héritage de priorité. Par exemple, si une tâche moins prioritaire acquiert une ressource et qu’ensuite une tâche plus prioritaire démarre et prend une autre ressource, puis ensuite nécessite la première ressource. Elle transfert alors sa priorité à l’autre tâche. Cependant, si cette autre tâche a aussi besoin de cette seconde ressource avant de terminer avec la première, les deux tâches deviennent bloquées. Dans le cas du ICCP, la tâche moins prioritaire obtiendrait déjà une plus grande priorité, jusqu’à la fin de l’utilisation des ressources dont la tâche plus prioritaire pourrait avoir besoin.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng1.InsertXML($xml1)

# --- Edit 2: the "_GoBack" bookmark used to sit at the very end of the document;
# now that the last real edit is in the new third paragraph above, move it there
# (already added by Edit 1) and strip the stale pair from the final paragraph.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rngLast = $pLast.Range

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="004926C7" w:rsidRDefault="004926C7" w:rsidP="004926C7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Le laboratoire était très bien structuré. Le fait d’avoir un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pdf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pour les questions et un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pdf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pour résumer la matière utile au laboratoire était très apprécié. Cela nous a permis de plus nous concentrer sur le laboratoire et moins sur la recherches d’informations dans nos notes de cours. Le temps consacré était approprié. Seulement quelques heures en dehors des heures de laboratoires fut nécessaires.</w:t></w:r><w:r w:rsidR="00346D9E"><w:t xml:space="preserve"> Somme toute, ce laboratoire était un excellent rappel des principes de bases du cours de noyau de système d’exploitation. </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rngLast.InsertXML($xml2)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
